$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current values of the columns that move (A,B,D,E,F,G,H,Q,R,Z,AB)
# for rows 3-6, then write them back rotated: row 6 -> row 3, and rows 3,4,5
# each shift down one row (3->4, 4->5, 5->6).

$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")
$rows = @(3,4,5,6)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# New row order: row3 gets old row6, row4 gets old row3, row5 gets old row4, row6 gets old row5
$mapping = @{ 3 = 6; 4 = 3; 5 = 4; 6 = 5 }

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcData[$c]
    }
}
